# praks registreerimine faili kirjutamine
# Adds a new "Praktikumid ja tunnivälised kursused" column (J) and
# reorders/updates the student rows 2-7 to match the refreshed roster.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column J header
$ws.Cells.Item(1, 10).Value = "Praktikumid ja tunnivälised kursused"

# Row 2 - Sarah 0 (columns A-I unchanged; J stays empty but present)
$ws.Cells.Item(2, 1).Value = "Sarah 0"
$ws.Cells.Item(2, 2).Value = "Keskkonnakeemia 2. periood"
$ws.Cells.Item(2, 3).Value = "Hispaania keel 1"
$ws.Cells.Item(2, 4).Value = "Filosoofia 1"
$ws.Cells.Item(2, 5).Value = "Hispaania keel 2"
$ws.Cells.Item(2, 6).Value = "Matemaatika ajaloo elemente ja rakendusi"
$ws.Cells.Item(2, 7).Value = "Hispaania keel 3"
$ws.Cells.Item(2, 8).Value = "Majandusmatemaatika elemendid"
$ws.Cells.Item(2, 9).Value = "Projektikirjutamisõpe"
$ws.Cells.Item(2, 10).NumberFormat = "General"
$ws.Cells.Item(2, 10).Value = ""

# Row 3 - Brianna 1
$ws.Cells.Item(3, 1).Value = "Brianna 1"
$ws.Cells.Item(3, 2).Value = "Joonestamine 2. periood"
$ws.Cells.Item(3, 3).Value = "Linux Raspberry Pi näitel"
$ws.Cells.Item(3, 4).Value = "Keemilised elemendid"
$ws.Cells.Item(3, 5).Value = "Keskkonnakeemia 3. periood"
$ws.Cells.Item(3, 6).Value = "Loogika 4. periood"
$ws.Cells.Item(3, 7).Value = "Teater Vanemuise kultuuritänavas"
$ws.Cells.Item(3, 8).Value = "Ettevõtlusõpe 5. periood"
$ws.Cells.Item(3, 9).Value = "CAD joonestamine"
$ws.Cells.Item(3, 10).Value = "Koorilaul"

# Row 4 - Lauren 4
$ws.Cells.Item(4, 1).Value = "Lauren 4"
$ws.Cells.Item(4, 2).Value = "Finantsmõtlemine"
$ws.Cells.Item(4, 3).Value = "Labortöid füüsikas 10. ja 11. klassile"
$ws.Cells.Item(4, 4).Value = "Ajakirjanduse alused"
$ws.Cells.Item(4, 5).Value = "Statistiline maailmapilt"
$ws.Cells.Item(4, 6).Value = "Matemaatika ajaloo elemente ja rakendusi"
$ws.Cells.Item(4, 7).Value = "Majandusõpe"
$ws.Cells.Item(4, 8).Value = "Laboratoorsed tööd bioloogias"
$ws.Cells.Item(4, 9).Value = "Mobiilirakenduste loomine (APP Inventor)"
$ws.Cells.Item(4, 10).Value = "Akvaristika, Koorilaul, Rahvatants"

# Row 5 - Hannah 3
$ws.Cells.Item(5, 1).Value = "Hannah 3"
$ws.Cells.Item(5, 2).Value = "Geoinformaatika"
$ws.Cells.Item(5, 3).Value = "Linux Raspberry Pi näitel"
$ws.Cells.Item(5, 4).Value = "Loomade käitumine 3. periood"
$ws.Cells.Item(5, 5).Value = "Turundus"
$ws.Cells.Item(5, 6).Value = "Köögifüüsika"
$ws.Cells.Item(5, 7).Value = "Millest ELU koosneb?"
$ws.Cells.Item(5, 8).Value = "Tänavakunst"
$ws.Cells.Item(5, 9).Value = "Projektikirjutamisõpe"
$ws.Cells.Item(5, 10).Value = "Koorilaul, Rahvatants"

# Row 6 - Jessica 2
$ws.Cells.Item(6, 1).Value = "Jessica 2"
$ws.Cells.Item(6, 2).Value = "Globaliseeruv maailm"
$ws.Cells.Item(6, 3).Value = "Õpioskused"
$ws.Cells.Item(6, 4).Value = "Küberkaitse 1"
$ws.Cells.Item(6, 5).Value = "Keskkonnakeemia 3. periood"
$ws.Cells.Item(6, 6).Value = "Loomade käitumine 4. periood"
$ws.Cells.Item(6, 7).Value = "Millest ELU koosneb?"
$ws.Cells.Item(6, 8).Value = "Tänavakunst"
$ws.Cells.Item(6, 9).Value = "Mobiilirakenduste loomine (APP Inventor)"
$ws.Cells.Item(6, 10).Value = "Rahvatants"

# Row 7 - Emma 5
$ws.Cells.Item(7, 1).Value = "Emma 5"
$ws.Cells.Item(7, 2).Value = "Geoinformaatika"
$ws.Cells.Item(7, 3).Value = "Inimene ja õigus"
$ws.Cells.Item(7, 4).Value = "Muusikaline kirjaoskus"
$ws.Cells.Item(7, 5).Value = "Inimene ja ühiskond 3. periood"
$ws.Cells.Item(7, 6).Value = "Joonestamine 4. periood"
$ws.Cells.Item(7, 7).Value = "Teater Vanemuise kultuuritänavas"
$ws.Cells.Item(7, 8).Value = "Ettevõtlusõpe 5. periood"
$ws.Cells.Item(7, 9).Value = "Karjääriõpetus"
$ws.Cells.Item(7, 10).Value = "Rahvatants, Koorilaul, Näitering, Akvaristika"
